$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 43317
$ws.Range("A2").Copy()
$ws.Range("A6").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("B6").Value = 0.5
$ws.Range("B2").Copy()
$ws.Range("B6").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("D6").Value = "setting up github"

$ws.Range("G8").Select()
